$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C ("Förändrad") from 45205 to 45206 for every existing
#    data row (rows 2 through 66).
for ($r = 2; $r -le 66; $r++) {
    $ws.Cells.Item($r, 3).Value = 45206
}

# 2) Ensure row 66 carries an explicit row height (matches the other data
#    rows, which already have ht="15" customHeight="1").
$ws.Rows.Item(66).RowHeight = 15

# 3) Append the new row 67 with the new logging-notification record.
$ws.Cells.Item(67, 1).Value = "A 48245-2023"
$ws.Cells.Item(67, 2).Value = 45205
$ws.Cells.Item(67, 3).Value = 45206
$ws.Cells.Item(67, 4).Value = "VÄSTMANLANDS LÄN"
$ws.Cells.Item(67, 5).Value = "HALLSTAHAMMAR"
$ws.Cells.Item(67, 6).Value = "Allmännings- och besparingsskogar"
$ws.Cells.Item(67, 7).Value = 15.7
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 13).Value = 0
$ws.Cells.Item(67, 14).Value = 0
$ws.Cells.Item(67, 15).Value = 0
$ws.Cells.Item(67, 16).Value = 0
$ws.Cells.Item(67, 17).Value = 0

# Match the date-formatted style used by columns B/C on the existing rows.
$ws.Cells.Item(67, 2).NumberFormat = $ws.Cells.Item(66, 2).NumberFormat
$ws.Cells.Item(67, 3).NumberFormat = $ws.Cells.Item(66, 3).NumberFormat

# Column R stays an (empty) wrap-text styled cell like the other recent rows.
$ws.Cells.Item(67, 18).WrapText = $true
